$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Points = 2)
$ws.Range("B3").Value = 2.13937
$ws.Range("C3").Value = 0.93955
$ws.Range("D3").Value = -10.25588

# Row 4 (Points = 3)
$ws.Range("B4").Value = -0.466
$ws.Range("C4").Value = -0.19561
$ws.Range("D4").Value = -1.32437

# Row 5 (Points = 4)
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = -8.62524
$ws.Range("D5").Value = 2.21758

# Row 6 (Points = 5)
$ws.Range("B6").Value = -0.04764
$ws.Range("C6").Value = 0.9313
$ws.Range("D6").Value = 7.91342
$ws.Range("E6").Value = "Fail"

# Row 7 (Points = 6)
$ws.Range("B7").Value = -1.37523
$ws.Range("C7").Value = 1.30298
$ws.Range("D7").Value = 1.35261
$ws.Range("E7").Value = "Pass"

# Row 8 (Points = 7)
$ws.Range("B8").Value = -1.36607
$ws.Range("C8").Value = -0.89624
$ws.Range("D8").Value = -0.52759
$ws.Range("E8").Value = "Pass"

# Row 9 (Points = 8)
$ws.Range("B9").Value = -1.66409
$ws.Range("C9").Value = -0.32369
$ws.Range("D9").Value = 0.68336
$ws.Range("E9").Value = "Pass"

# Row 10 (Points = 9)
$ws.Range("B10").Value = -1.18109
$ws.Range("C10").Value = 0.01451
$ws.Range("D10").Value = -3.23807
$ws.Range("E10").Value = "Fail"

# Row 11 (Points = 10)
$ws.Range("B11").Value = -0.75213
$ws.Range("C11").Value = 0.8610100000000001
$ws.Range("D11").Value = 1.42462
